$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6541
$ws.Range("C27").Value = 1018
$ws.Range("D27").Value = 6099524
$ws.Range("E27").Value = 932.5063445956275
$ws.Range("F27").Value = 9.932773109243698
$ws.Range("G27").Value = 7.4973600844773
$ws.Range("H27").Value = 25.10543006091468
